$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F5, F7, F11, F12, F13
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 2807
$ws1.Range("F7").Value = 228
$ws1.Range("F11").Value = 85
$ws1.Range("F12").Value = 2636
$ws1.Range("F13").Value = 821

# Sheet "全部类型" (sheet4): update F6, F8, F13, F14, F15
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 2807
$ws4.Range("F8").Value = 228
$ws4.Range("F13").Value = 85
$ws4.Range("F14").Value = 2636
$ws4.Range("F15").Value = 821
